$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing last row (307): high/close were revised ---
$ws.Range("D307").Value = 460.02
$ws.Range("F307").Value = 459.95

# --- Append 3 new monthly rows (308-310), matching the style of row 307 ---
# Copy formatting (date style) from A307 down into the new A308:A310 cells first
$ws.Range("A307").Copy($ws.Range("A308:A310"))

# Row 308
$ws.Range("A308").Value = 45047.33333333334
$ws.Range("B308").Value = "FX_IDC:USDNGN"
$ws.Range("C308").Value = 459.82
$ws.Range("D308").Value = 464.98
$ws.Range("E308").Value = 459.82
$ws.Range("F308").Value = 460.58
$ws.Range("G308").Value = 0

# Row 309
$ws.Range("A309").Value = 45078.33333333334
$ws.Range("B309").Value = "FX_IDC:USDNGN"
$ws.Range("C309").Value = 460.58
$ws.Range("D309").Value = 820
$ws.Range("E309").Value = 460.58
$ws.Range("F309").Value = 758.77
$ws.Range("G309").Value = 0

# Row 310
$ws.Range("A310").Value = 45110.33333333334
$ws.Range("B310").Value = "FX_IDC:USDNGN"
$ws.Range("C310").Value = 758.77
$ws.Range("D310").Value = 774
$ws.Range("E310").Value = 758.77
$ws.Range("F310").Value = 767
$ws.Range("G310").Value = 0
